$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 98
$ws.Range("H98").Value = 364243.9
$ws.Range("I98").Value = 7752.96
$ws.Range("J98").Value = 3335001.8
$ws.Range("K98").Value = 7752.96
$ws.Range("L98").Value = 3335001.8
$ws.Range("M98").Value = -6254.96
$ws.Range("N98").Value = -3337997.8

# Row 122
$ws.Range("H122").Value = 364243.9
$ws.Range("I122").Value = 7752.96
$ws.Range("J122").Value = 3335001.8
$ws.Range("K122").Value = 23258.88
$ws.Range("L122").Value = 10005005.4
$ws.Range("M122").Value = -20808.88
$ws.Range("N122").Value = -10009905.4

# Row 137
$ws.Range("H137").Value = 789.95654
$ws.Range("I137").Value = 685.7826
$ws.Range("J137").Value = 894.13043
$ws.Range("K137").Value = 2057.3478
$ws.Range("L137").Value = 2682.39129
$ws.Range("M137").Value = 492.6522
$ws.Range("N137").Value = -7782.39129

# Row 138
$ws.Range("H138").Value = 3386.6616
$ws.Range("I138").Value = 1919.8108
$ws.Range("J138").Value = 5325
$ws.Range("K138").Value = 5759.4324
$ws.Range("L138").Value = 15975
$ws.Range("M138").Value = -619.4323999999997
$ws.Range("N138").Value = -26255


# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 74
$ws.Range("H74").Value = 962.3617
$ws.Range("I74").Value = 929.29266
$ws.Range("J74").Value = 1188.3334
$ws.Range("K74").Value = 929.29266
$ws.Range("L74").Value = 1188.3334
$ws.Range("M74").Value = -55.29265999999996
$ws.Range("N74").Value = -2936.3334

# Row 77
$ws.Range("H77").Value = 962.3617
$ws.Range("I77").Value = 929.29266
$ws.Range("J77").Value = 1188.3334
$ws.Range("K77").Value = 4646.463299999999
$ws.Range("L77").Value = 5941.666999999999
$ws.Range("M77").Value = -278.4632999999994
$ws.Range("N77").Value = -14677.667

# Row 102
$ws.Range("H102").Value = 4702.8184
$ws.Range("I102").Value = 3215
$ws.Range("K102").Value = 3215
$ws.Range("M102").Value = -1593


# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 94
$ws.Range("H94").Value = 838.6667
$ws.Range("I94").Value = 825.28
$ws.Range("J94").Value = 869.0909
$ws.Range("K94").Value = 825.28
$ws.Range("L94").Value = 869.0909
$ws.Range("M94").Value = -374.28
$ws.Range("N94").Value = -1771.0909

# Row 99
$ws.Range("H99").Value = 2841.2
$ws.Range("I99").Value = 1080.2
$ws.Range("J99").Value = 4602.2
$ws.Range("K99").Value = 1080.2
$ws.Range("L99").Value = 4602.2
$ws.Range("M99").Value = 417.8
$ws.Range("N99").Value = -7598.2

# Row 105
$ws.Range("H105").Value = 1552.7179
$ws.Range("I105").Value = 1449.8636
$ws.Range("K105").Value = 1449.8636
$ws.Range("M105").Value = 297.1364000000001

# Row 107
$ws.Range("H107").Value = 2467.25
$ws.Range("I107").Value = 1660.3334
$ws.Range("J107").Value = 4888
$ws.Range("K107").Value = 1660.3334
$ws.Range("L107").Value = 4888
$ws.Range("M107").Value = 259.6666
$ws.Range("N107").Value = -8728


# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 122
$ws.Range("H122").Value = 512
$ws.Range("I122").Value = 512
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 1536
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = 914
$ws.Range("N122").ClearContents()

# Row 132
$ws.Range("H132").Value = 1828.8966
$ws.Range("I132").Value = 1617.3914
$ws.Range("J132").Value = 2639.6667
$ws.Range("K132").Value = 4852.174199999999
$ws.Range("L132").Value = 7919.000100000001
$ws.Range("M132").Value = -2322.174199999999
$ws.Range("N132").Value = -12979.0001


# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 5
$ws.Range("H5").Value = 711
$ws.Range("I5").Value = 456.2857
$ws.Range("K5").Value = 1368.8571
$ws.Range("M5").Value = -1256.8571

# Row 122
$ws.Range("H122").Value = 567
$ws.Range("I122").Value = 319.54544
$ws.Range("J122").Value = 814.4545000000001
$ws.Range("K122").Value = 2875.90896
$ws.Range("L122").Value = 7330.0905
$ws.Range("M122").Value = -425.9089599999998
$ws.Range("N122").Value = -12230.0905

# Row 131
$ws.Range("H131").Value = 6850167.5
$ws.Range("J131").Value = 8334191
$ws.Range("L131").Value = 25002573
$ws.Range("N131").Value = -25012653

# Row 135
$ws.Range("H135").Value = 711
$ws.Range("I135").Value = 456.2857
$ws.Range("K135").Value = 4106.571300000001
$ws.Range("M135").Value = -1571.571300000001


# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 10000
$ws.Range("I15").Value = 0
$ws.Range("J15").Value = 10000
$ws.Range("K15").Value = 0
$ws.Range("L15").Value = 10000
$ws.Range("M15").ClearContents()
$ws.Range("N15").Value = -10576

# Row 74
$ws.Range("H74").Value = 0
$ws.Range("J74").Value = 0
$ws.Range("L74").Value = 0
$ws.Range("N74").ClearContents()

# Row 77
$ws.Range("H77").Value = 0
$ws.Range("J77").Value = 0
$ws.Range("L77").Value = 0
$ws.Range("N77").ClearContents()

# Row 81
$ws.Range("H81").Value = 10000
$ws.Range("I81").Value = 0
$ws.Range("J81").Value = 10000
$ws.Range("K81").Value = 0
$ws.Range("L81").Value = 10000
$ws.Range("M81").ClearContents()
$ws.Range("N81").Value = -11996

# Row 84
$ws.Range("H84").Value = 10000
$ws.Range("I84").Value = 0
$ws.Range("J84").Value = 10000
$ws.Range("K84").Value = 0
$ws.Range("L84").Value = 30000
$ws.Range("M84").ClearContents()
$ws.Range("N84").Value = -39984

# Row 123
$ws.Range("H123").Value = 17559.5
$ws.Range("J123").Value = 17559.5
$ws.Range("L123").Value = 17559.5
$ws.Range("N123").Value = -22459.5

# Row 132
$ws.Range("H132").Value = 2235.6785
$ws.Range("I132").Value = 2285
$ws.Range("J132").Value = 2087.7144
$ws.Range("K132").Value = 6855
$ws.Range("L132").Value = 6263.1432
$ws.Range("M132").Value = -4325
$ws.Range("N132").Value = -11323.1432


# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 81
$ws.Range("H81").Value = 2383.9443
$ws.Range("I81").Value = 1562.3846
$ws.Range("J81").Value = 4520
$ws.Range("K81").Value = 3124.7692
$ws.Range("L81").Value = 9040
$ws.Range("M81").Value = -2063.7692
$ws.Range("N81").Value = -11162

# Row 82
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("N82").ClearContents()

# Row 84
$ws.Range("H84").Value = 2383.9443
$ws.Range("I84").Value = 1562.3846
$ws.Range("J84").Value = 4520
$ws.Range("K84").Value = 15623.846
$ws.Range("L84").Value = 45200
$ws.Range("M84").Value = -10319.846
$ws.Range("N84").Value = -55808

# Row 85
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("N85").ClearContents()

# Row 122
$ws.Range("H122").Value = 1496.2693
$ws.Range("I122").Value = 1451
$ws.Range("J122").Value = 1568.7
$ws.Range("K122").Value = 4353
$ws.Range("L122").Value = 4706.1
$ws.Range("M122").Value = -1903
$ws.Range("N122").Value = -9606.1

# Row 132
$ws.Range("H132").Value = 20162230
$ws.Range("I132").Value = 24510682
$ws.Range("J132").Value = 1228.1818
$ws.Range("K132").Value = 73532046
$ws.Range("L132").Value = 3684.5454
$ws.Range("M132").Value = -73529516
$ws.Range("N132").Value = -8744.545399999999

# Row 136
$ws.Range("H136").Value = 582.4167
$ws.Range("I136").Value = 500.35483
$ws.Range("J136").Value = 732.05884
$ws.Range("K136").Value = 1501.06449
$ws.Range("L136").Value = 2196.17652
$ws.Range("M136").Value = 1048.93551
$ws.Range("N136").Value = -7296.17652

